$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'268.21"
$ws.Range("D3").Value = "'22.89"
$ws.Range("D5").Value = "'0.06186"
$ws.Range("D6").Value = "'3.583"
$ws.Range("D7").Value = "'6.699"
$ws.Range("D8").Value = "'1.367"
$ws.Range("D9").Value = "'0.8409"
$ws.Range("D10").Value = "'0.01364"
$ws.Range("D12").Value = "'0.08239"
$ws.Range("D13").Value = "'0.03421"
$ws.Range("D14").Value = "'0.03267"
$ws.Range("D15").Value = "'0.09274"
$ws.Range("D16").Value = "'3.907"
$ws.Range("D17").Value = "'0.001733"
$ws.Range("D18").Value = "'0.04861"
$ws.Range("D19").Value = "'0.006258"
$ws.Range("D20").Value = "'0.005359"
$ws.Range("D21").Value = "'0.001090"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("D23").Value = "'3.757"
$ws.Range("D24").Value = "'2.352"
$ws.Range("D26").Value = "'0.1213"
$ws.Range("D40").Value = "'0.04675"
$ws.Range("D41").Value = "'0.006958"
$ws.Range("D42").Value = "'0.1153"
$ws.Range("D43").Value = "'0.003351"
$ws.Range("D44").Value = "'0.01216"
$ws.Range("D45").Value = "'0.00006242"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D47").Value = "'0.7005"
$ws.Range("D48").Value = "'0.1620"
$ws.Range("D49").Value = "'0.00002102"
